$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 2044
$ws.Range("I101").Value = 232.66667
$ws.Range("J101").Value = 3855.3333
$ws.Range("K101").Value = 698.00001
$ws.Range("L101").Value = 11565.9999
$ws.Range("M101").Value = 923.99999
$ws.Range("N101").Value = -14809.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3033.0952
$ws.Range("I137").Value = 2946.8
$ws.Range("J137").Value = 3248.8333
$ws.Range("K137").Value = 8840.400000000001
$ws.Range("L137").Value = 9746.499899999999
$ws.Range("M137").Value = -6290.400000000001
$ws.Range("N137").Value = -14846.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2673.2856
$ws.Range("I45").Value = 1825
$ws.Range("J45").Value = 3444.4546
$ws.Range("K45").Value = 1825
$ws.Range("L45").Value = 3444.4546
$ws.Range("M45").Value = -1448
$ws.Range("N45").Value = -4198.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8463824
$ws.Range("I61").Value = 10528613
$ws.Range("K61").Value = 10528613
$ws.Range("M61").Value = -10528401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 957.6486
$ws.Range("I97").Value = 776.5862
$ws.Range("J97").Value = 1614
$ws.Range("K97").Value = 776.5862
$ws.Range("L97").Value = 1614
$ws.Range("M97").Value = -280.5862
$ws.Range("N97").Value = -2606

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4196.1924
$ws.Range("I122").Value = 4272.5264
$ws.Range("K122").Value = 12817.5792
$ws.Range("M122").Value = -10367.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8463824
$ws.Range("I136").Value = 10528613
$ws.Range("K136").Value = 31585839
$ws.Range("M136").Value = -31583289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 760
$ws.Range("I8").Value = 1100
$ws.Range("J8").Value = 533.3333
$ws.Range("K8").Value = 1100
$ws.Range("L8").Value = 533.3333
$ws.Range("M8").Value = -960
$ws.Range("N8").Value = -813.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1809.7106
$ws.Range("I94").Value = 1963.8334
$ws.Range("J94").Value = 1231.75
$ws.Range("K94").Value = 1963.8334
$ws.Range("L94").Value = 1231.75
$ws.Range("M94").Value = -1512.8334
$ws.Range("N94").Value = -2133.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 99999
$ws.Range("J70").Value = 99999
$ws.Range("L70").Value = 99999
$ws.Range("N70").Value = -100629

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 99999
$ws.Range("J73").Value = 99999
$ws.Range("L73").Value = 99999
$ws.Range("N73").Value = -102183

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2800.9333
$ws.Range("J12").Value = 3946.125
$ws.Range("L12").Value = 11838.375
$ws.Range("N12").Value = -12184.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1009143
$ws.Range("I32").Value = 2501595
$ws.Range("J32").Value = 14175
$ws.Range("K32").Value = 7504785
$ws.Range("L32").Value = 42525
$ws.Range("M32").Value = -7504502
$ws.Range("N32").Value = -43091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 7445.857
$ws.Range("I44").Value = 429.5
$ws.Range("K44").Value = 1288.5
$ws.Range("M44").Value = -890.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 10264.223
$ws.Range("J74").Value = 13118.286
$ws.Range("L74").Value = 39354.858
$ws.Range("N74").Value = -41476.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 10264.223
$ws.Range("J77").Value = 13118.286
$ws.Range("L77").Value = 118064.574
$ws.Range("N77").Value = -128672.574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 5531.8
$ws.Range("I108").Value = 2442.7778
$ws.Range("K108").Value = 7328.3334
$ws.Range("M108").Value = -4448.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 6845.2856
$ws.Range("I109").Value = 916.8
$ws.Range("K109").Value = 2750.4
$ws.Range("M109").Value = -1710.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 12888.167
$ws.Range("I118").Value = 8799.200000000001
$ws.Range("K118").Value = 26397.6
$ws.Range("M118").Value = -25154.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 7952.0835
$ws.Range("I136").Value = 2870.7144
$ws.Range("K136").Value = 8612.143199999999
$ws.Range("M136").Value = -3512.143199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1935.8778
$ws.Range("I139").Value = 1939.2941
$ws.Range("K139").Value = 5817.8823
$ws.Range("M139").Value = -677.8823000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2053.6667
$ws.Range("J80").Value = 2504
$ws.Range("L80").Value = 2504
$ws.Range("N80").Value = -4500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2053.6667
$ws.Range("J83").Value = 2504
$ws.Range("L83").Value = 12520
$ws.Range("N83").Value = -22504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2182.3684
$ws.Range("I102").Value = 2059.8125
$ws.Range("J102").Value = 2836
$ws.Range("K102").Value = 2059.8125
$ws.Range("L102").Value = 2836
$ws.Range("M102").Value = -437.8125
$ws.Range("N102").Value = -6080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3286.75
$ws.Range("I122").Value = 2866.5557
$ws.Range("J122").Value = 3827
$ws.Range("K122").Value = 8599.667099999999
$ws.Range("L122").Value = 11481
$ws.Range("M122").Value = -6149.667099999999
$ws.Range("N122").Value = -16381

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7588.647
$ws.Range("I126").Value = 8687.929
$ws.Range("K126").Value = 26063.787
$ws.Range("M126").Value = -23593.787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3450714
$ws.Range("J132").Value = 14289091
$ws.Range("L132").Value = 42867273
$ws.Range("N132").Value = -42872333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 68204
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4145.1665
$ws.Range("I7").Value = 4145.1665
$ws.Range("K7").Value = 4145.1665
$ws.Range("M7").Value = -4033.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2623.2354
$ws.Range("I40").Value = 2439.8667
$ws.Range("J40").Value = 3998.5
$ws.Range("K40").Value = 2439.8667
$ws.Range("L40").Value = 3998.5
$ws.Range("M40").Value = -2303.8667
$ws.Range("N40").Value = -4270.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3326.074
$ws.Range("I122").Value = 2897.818
$ws.Range("K122").Value = 8693.454000000002
$ws.Range("M122").Value = -6243.454000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4145.1665
$ws.Range("I126").Value = 4145.1665
$ws.Range("K126").Value = 12435.4995
$ws.Range("M126").Value = -9965.499500000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4191.7144
$ws.Range("I132").Value = 2982.182
$ws.Range("J132").Value = 5522.2
$ws.Range("K132").Value = 8946.545999999998
$ws.Range("L132").Value = 16566.6
$ws.Range("M132").Value = -6416.545999999998
$ws.Range("N132").Value = -21626.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 55200
$ws.Range("J134").Value = 55200
$ws.Range("L134").Value = 55200
$ws.Range("N134").Value = -65340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 5002.5
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1348.4
$ws.Range("I81").Value = 1348.4
$ws.Range("K81").Value = 2696.8
$ws.Range("M81").Value = -1635.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1348.4
$ws.Range("I84").Value = 1348.4
$ws.Range("K84").Value = 13484
$ws.Range("M84").Value = -8180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2977.276
$ws.Range("I122").Value = 2804.48
$ws.Range("J122").Value = 4057.25
$ws.Range("K122").Value = 8413.440000000001
$ws.Range("L122").Value = 12171.75
$ws.Range("M122").Value = -5963.440000000001
$ws.Range("N122").Value = -17071.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3463.8965
$ws.Range("I126").Value = 3289.739
$ws.Range("K126").Value = 9869.217000000001
$ws.Range("M126").Value = -7399.217000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 240241.4
$ws.Range("I132").Value = 1969.1786
$ws.Range("K132").Value = 5907.5358
$ws.Range("M132").Value = -3377.5358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 30657
$ws.Range("J133").Value = 30657
$ws.Range("L133").Value = 30657
$ws.Range("N133").Value = -40777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 127108.17
$ws.Range("I136").Value = 7089.388
$ws.Range("J136").Value = 629686.8
$ws.Range("K136").Value = 21268.164
$ws.Range("L136").Value = 1889060.4
$ws.Range("M136").Value = -18718.164
$ws.Range("N136").Value = -1894160.4
